# Covid19.xlsx update: append new daily rows (40-45) of COVID data to the
# Summary sheet, refresh the Death-Projection what-if inputs, and nudge a
# couple of chart / drawing layouts to match the re-saved workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet: new raw data for rows 40-45 (C/D/E) + formulas that
#    mirror the pattern already used by the rows above (F..R, T).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$newRows = @(
    @{ Row = 40; C = 461437; D = 16478; E = 25410 },
    @{ Row = 41; C = 496535; D = 18586; E = 28790 },
    @{ Row = 42; C = 526396; D = 20463; E = 31270 },
    @{ Row = 43; C = 555313; D = 22020; E = 32988 },
    @{ Row = 44; C = 580619; D = 23529; E = 43482 },
    @{ Row = 45; C = 607670; D = 25832; E = 47763 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prev = $row - 1

    $summary.Range("C$row").Value = $r.C
    $summary.Range("D$row").Value = $r.D
    $summary.Range("E$row").Value = $r.E

    $summary.Range("F$row").Formula = "=D$row*`$W`$7"
    $summary.Range("G$row").Formula = "=E$row*`$W`$7"
    $summary.Range("H$row").Formula = "=C$row-C$prev"
    $summary.Range("I$row").Formula = "=D$row-D$prev"
    $summary.Range("J$row").Formula = "=J$prev+1"
    $summary.Range("K$row").Formula = "=W`$4*EXP(X`$4*J$row)"
    $summary.Range("L$row").Value = "log"
    $summary.Range("M$row").Formula = "=SLOPE(`$C`$34:C$row,LN(`$J`$34:J$row))"
    $summary.Range("N$row").Formula = "=INTERCEPT(`$C`$34:C$row, LN(`$J`$34:J$row))"
    $summary.Range("O$row").Formula = "=INDEX(LINEST(`$C`$28:C$row,LN(`$J`$28:J$row),1,1),3)"
    $summary.Range("P$row").Formula = "=M$row*LN(J$row)+N$row"
    $summary.Range("Q$row").Formula = "=P$row-C$row"
    $summary.Range("R$row").Formula = "=24*(T$row-J$row)"
    $summary.Range("T$row").Formula = "=EXP((2*M$row*LN(J$row)+N$row)/M$row)"
}

# frozen pane / selection follow the new bottom of the data block
$summary.Activate()
$summary.Range("M45").Select()

# ---------------------------------------------------------------------
# 2. Death Projection sheet: refresh the "today" inputs used by the
#    what-if block (doubling window shrinks from 9 to 5 days, and the
#    lookup date moves back 10 days).
# ---------------------------------------------------------------------
$deathProjection = $wb.Worksheets.Item("Death Projection")
$deathProjection.Range("D5").Value = 5
$deathProjection.Range("B6").Value = 43918
$deathProjection.Activate()
$deathProjection.Range("D13").Select()

# ---------------------------------------------------------------------
# 3. Charts sheet: keep it the active tab, matching the saved selection.
# ---------------------------------------------------------------------
$charts = $wb.Worksheets.Item("Charts")
$charts.Activate()
$charts.Range("H24").Select()

$wb.Save()
